# Updates cryptos list (price / 1h-volume-change columns, plus two swapped
# rows for Cardano/OKB) to match the "Thu Nov 9 03:14:33 UTC 2023" refresh.
#
# Note: several "Price" (column D) values are numeric-looking strings
# (e.g. "248.87", "1.00") that Excel's COM layer would otherwise coerce to
# real numbers on plain assignment, silently dropping the original text
# formatting (trailing zeros, etc.). For those cells we briefly force a
# text NumberFormat before writing the value, then restore the cell style
# to "Normal" so the cell's style index/appearance is unchanged afterward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.328.03'
$ws.Range("E2").Value = '  +2.84%  '
$ws.Range("D3").Value = '1.915.95'
$ws.Range("E3").Value = '  +1.59%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("E6").Value = '  +0.51%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.88'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.80%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '57.71'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.46%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.365'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0761'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0993'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.39'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.798'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.72%  '
$ws.Range("D15").Value = '2.198.42'
$ws.Range("E15").Value = '  +1.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.59%  '
$ws.Range("D17").Value = '1.935.76'
$ws.Range("E17").Value = '  +2.34%  '
$ws.Range("D18").Value = '36.309.24'
$ws.Range("E18").Value = '  +2.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.71%  '
$ws.Range("D20").Value = '0.0₃0845'
$ws.Range("E20").Value = '  +2.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '251.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.70'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '167.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.77'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.80'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.81%  '
$ws.Range("E30").Value = '  +1.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.53'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0607'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.95'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.31'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.58%  '
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("E36").Value = '  +21.34%  '
$ws.Range("E37").Value = '  -14.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.858'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.52%  '
$ws.Range("E39").Value = '  +2.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '104.85'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0228'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.28'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +24.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.05'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.50%  '
$ws.Range("E44").Value = '  +2.87%  '
$ws.Range("D45").Value = '1.340.36'
$ws.Range("E45").Value = '  +3.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0807'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.46%  '
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("E49").Value = '  +1.91%  '
$ws.Range("E50").Value = '  +3.07%  '
$ws.Range("D51").Value = '2.097.33'
$ws.Range("E51").Value = '  +1.46%  '
